$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of student data below the existing sample row.
# A3/B3/C3/F3 pick up a "pasted text" style (numFmtId 49, distinct font),
# matching how these four cells were filled in from clipboard data, while
# D3/E3 were left with their plain column-default formatting.
$ws.Range("A3").Value = "999"
$ws.Range("A3").Font.Family = 4

$ws.Range("B3").Value = "测试"
$ws.Range("B3").Font.Family = 4

$ws.Range("C3").Value = "男"
$ws.Range("C3").Font.Family = 4

$ws.Range("D3").Value = 2021

$ws.Range("E3").Value = "高二"

$ws.Range("F3").Value = "高二2班"
$ws.Range("F3").Font.Family = 4

$ws.Range("F8").Select()
